# Update column F ("dSF") values for specific rows as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    5  = -2
    6  = 0
    13 = -3
    14 = 1
    16 = 6
    25 = 2
    28 = 1
    34 = 1
    37 = 0
    40 = -1
    41 = -1
    44 = -3
    45 = -1
    46 = 6
    49 = -7
    53 = -2
    55 = -3
    63 = -8
    67 = -4
    68 = -8
    69 = -6
    71 = 1
    74 = 1
    81 = 1
    83 = -3
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
